$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the JSON, Common Built-in Methods and Fetch API rows as completed
$ws.Range("E26").Value = 1
$ws.Range("E27").Value = 1
$ws.Range("E28").Value = 1

# Move the active selection down to where the author left off editing
$ws.Range("E29").Select()
